$d = $word.ActiveDocument

# --- Change 1: bold "id" -> bold "I" + bold "d" (two separate runs) ---
$p1 = $d.Paragraphs(3)
$r1 = $p1.Range.Duplicate
$null = $r1.Find.Execute("id", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Text = "Id"
# Toggle Bold off/on for the second letter only, which forces Word to keep it
# as a distinct run (identical bold formatting) instead of re-merging it with
# the first letter's run.
$rSecond = $d.Range($r1.Start + 1, $r1.Start + 2)
$rSecond.Font.Bold = $false
$rSecond.Font.Bold = $true

# --- Change 2: IsHardCover description split into two runs with new wording ---
$p2 = $d.Paragraphs(15)
$r2 = $p2.Range.Duplicate
$null = $r2.Find.Execute(" = definisce se è una copertina rigida", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insStart = $r2.Start
# Clear the matched text completely, then insert the two pieces of text one
# after another into the now-empty range. Sequential InsertAfter calls on an
# originally empty range each create their own run without leftover
# formatting artifacts.
$r2.Text = ""
$firstPart = $d.Range($insStart, $insStart)
$firstPart.InsertAfter(" = definisce se è ")
$secondStart = $insStart + " = definisce se è ".Length
$secondPart = $d.Range($secondStart, $secondStart)
$secondPart.InsertAfter("rilegato")
